# Merge "Final Demand" related folders (SwFD, GDP projections, AIC) into a
# single "Final Demand" tree, and add the new "Projections" / "Merged FD"
# lookup rows, per commit "Tentativo di unire in unico excel la FD".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 15 (SwFD): carol/matti paths move from Recycling\SwFD to Final Demand\SwFD
$ws.Range("D15").Value = "C:\Users\carol\OneDrive\Documenti\GitHub\GreenTechs\Final Demand\SwFD"
$ws.Range("E15").Value = "C:\Users\matti\OneDrive - Politecnico di Milano\Documenti\GitHub\GreenTechs\Final Demand\SwFD"

# Row 18: "GDP projection" folder/label merges into "Final Demand" as "GDP projections"
$ws.Range("A18").Value = "GDP projections"
$ws.Range("D18").Value = "C:\Users\carol\OneDrive\Documenti\GitHub\GreenTechs\Final Demand\GDP projections.xlsx"
$ws.Range("E18").Value = "C:\Users\matti\OneDrive - Politecnico di Milano\Documenti\GitHub\GreenTechs\Final Demand\GDP projections.xlsx"

# Row 19 (AIC): carol path newly added, matti path moves from Recycling\AIC to Final Demand\AIC
$ws.Range("A19").Value = "AIC"
$ws.Range("D19").Value = "C:\Users\carol\OneDrive\Documenti\GitHub\GreenTechs\Final Demand\AIC"
$ws.Range("E19").Value = "C:\Users\matti\OneDrive - Politecnico di Milano\Documenti\GitHub\GreenTechs\Final Demand\AIC"

# Row 20 (new): "Projections" label pointing at the Final Demand folder
$ws.Range("A20").Value = "Projections"
$ws.Range("D20").Value = "C:\Users\carol\OneDrive\Documenti\GitHub\GreenTechs\Final Demand"

# Row 21 (new): "Merged FD" label pointing at the Final Demand folder
$ws.Range("A21").Value = "Merged FD"
$ws.Range("D21").Value = "C:\Users\carol\OneDrive\Documenti\GitHub\GreenTechs\Final Demand"

# Update the frozen-pane/selection view state to match the new extent.
$ws.Range("D5").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("D21").Select()
